$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1426182.8
$ws.Range("I11").Value = 1426182.8
$ws.Range("K11").Value = 1426182.8
$ws.Range("M11").Value = -1426042.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 999.6
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 998
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 998
$ws.Range("M43").Value = -931
$ws.Range("N43").Value = -1136

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1569.72
$ws.Range("I98").Value = 1081.4286
$ws.Range("J98").Value = 4133.25
$ws.Range("K98").Value = 1081.4286
$ws.Range("L98").Value = 4133.25
$ws.Range("M98").Value = 416.5714
$ws.Range("N98").Value = -7129.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2040.0667
$ws.Range("I116").Value = 1980.2
$ws.Range("K116").Value = 1980.2
$ws.Range("M116").Value = 1461.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1569.72
$ws.Range("I122").Value = 1081.4286
$ws.Range("J122").Value = 4133.25
$ws.Range("K122").Value = 3244.2858
$ws.Range("L122").Value = 12399.75
$ws.Range("M122").Value = -794.2857999999997
$ws.Range("N122").Value = -17299.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 49997.5
$ws.Range("J133").Value = 49997.5
$ws.Range("L133").Value = 49997.5
$ws.Range("N133").Value = -60117.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3038.5854
$ws.Range("J138").Value = 3780.76
$ws.Range("L138").Value = 11342.28
$ws.Range("N138").Value = -21622.28

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 4587.8
$ws.Range("I25").Value = 1984.75
$ws.Range("J25").Value = 15000
$ws.Range("K25").Value = 1984.75
$ws.Range("L25").Value = 15000
$ws.Range("M25").Value = -1582.75
$ws.Range("N25").Value = -15804

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1804.7059
$ws.Range("I74").Value = 1815.0322
$ws.Range("J74").Value = 1698
$ws.Range("K74").Value = 1815.0322
$ws.Range("L74").Value = 1698
$ws.Range("M74").Value = -941.0322000000001
$ws.Range("N74").Value = -3446

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1804.7059
$ws.Range("I77").Value = 1815.0322
$ws.Range("J77").Value = 1698
$ws.Range("K77").Value = 9075.161
$ws.Range("L77").Value = 8490
$ws.Range("M77").Value = -4707.161
$ws.Range("N77").Value = -17226

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3680.7896
$ws.Range("I20").Value = 4089.6155
$ws.Range("K20").Value = 4089.6155
$ws.Range("M20").Value = -3842.6155

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 843.7586
$ws.Range("I94").Value = 746.62964
$ws.Range("K94").Value = 746.62964
$ws.Range("M94").Value = -295.62964

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1835.2142
$ws.Range("I107").Value = 1707.1538
$ws.Range("J107").Value = 3500
$ws.Range("K107").Value = 1707.1538
$ws.Range("L107").Value = 3500
$ws.Range("M107").Value = 212.8462
$ws.Range("N107").Value = -7340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1940.8889
$ws.Range("I31").Value = 1223.5555
$ws.Range("J31").Value = 4810.222
$ws.Range("K31").Value = 1223.5555
$ws.Range("L31").Value = 4810.222
$ws.Range("M31").Value = -928.5554999999999
$ws.Range("N31").Value = -5400.222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1940.8889
$ws.Range("I34").Value = 1223.5555
$ws.Range("J34").Value = 4810.222
$ws.Range("K34").Value = 1223.5555
$ws.Range("L34").Value = 4810.222
$ws.Range("M34").Value = -1021.5555
$ws.Range("N34").Value = -5214.222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2787.3333
$ws.Range("I99").Value = 1712
$ws.Range("J99").Value = 3325
$ws.Range("K99").Value = 1712
$ws.Range("L99").Value = 3325
$ws.Range("M99").Value = -214
$ws.Range("N99").Value = -6321

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2787.3333
$ws.Range("I126").Value = 1712
$ws.Range("J126").Value = 3325
$ws.Range("K126").Value = 5136
$ws.Range("L126").Value = 9975
$ws.Range("M126").Value = -2666
$ws.Range("N126").Value = -14915

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 44.72222
$ws.Range("I2").Value = 27.8
$ws.Range("J2").Value = 47.451614
$ws.Range("K2").Value = 166.8
$ws.Range("L2").Value = 284.709684
$ws.Range("M2").Value = -53.80000000000001
$ws.Range("N2").Value = -510.709684

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1826.1154
$ws.Range("I34").Value = 349.42856
$ws.Range("J34").Value = 2370.158
$ws.Range("K34").Value = 1048.28568
$ws.Range("L34").Value = 7110.474
$ws.Range("M34").Value = -964.28568
$ws.Range("N34").Value = -7278.474

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2632
$ws.Range("I80").Value = 2800
$ws.Range("J80").Value = 2572
$ws.Range("K80").Value = 8400
$ws.Range("L80").Value = 7716
$ws.Range("M80").Value = -7464
$ws.Range("N80").Value = -9588

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 2632
$ws.Range("I83").Value = 2800
$ws.Range("J83").Value = 2572
$ws.Range("K83").Value = 25200
$ws.Range("L83").Value = 23148
$ws.Range("M83").Value = -20520
$ws.Range("N83").Value = -32508

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 749.75
$ws.Range("I86").Value = 649
$ws.Range("J86").Value = 783.3333
$ws.Range("K86").Value = 1947
$ws.Range("L86").Value = 2349.9999
$ws.Range("M86").Value = -761
$ws.Range("N86").Value = -4721.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 749.75
$ws.Range("I89").Value = 649
$ws.Range("J89").Value = 783.3333
$ws.Range("K89").Value = 5841
$ws.Range("L89").Value = 7049.9997
$ws.Range("M89").Value = 87
$ws.Range("N89").Value = -18905.9997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 646.34045
$ws.Range("I113").Value = 671.0263
$ws.Range("J113").Value = 542.1111
$ws.Range("K113").Value = 2013.0789
$ws.Range("L113").Value = 1626.3333
$ws.Range("M113").Value = 156.9211
$ws.Range("N113").Value = -5966.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1228.0625
$ws.Range("I132").Value = 834.9
$ws.Range("K132").Value = 7514.099999999999
$ws.Range("M132").Value = -4984.099999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 19150
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 19150
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 19150
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -19736

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2150.4
$ws.Range("I132").Value = 1958.2
$ws.Range("J132").Value = 2342.6
$ws.Range("K132").Value = 5874.6
$ws.Range("L132").Value = 7027.799999999999
$ws.Range("M132").Value = -3344.6
$ws.Range("N132").Value = -12087.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 938
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -205

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 938
$ws.Range("I27").Value = 500
$ws.Range("K27").Value = 500
$ws.Range("M27").Value = -393

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6354.8887
$ws.Range("I122").Value = 5833.0645
$ws.Range("J122").Value = 7510.357
$ws.Range("K122").Value = 17499.1935
$ws.Range("L122").Value = 22531.071
$ws.Range("M122").Value = -15049.1935
$ws.Range("N122").Value = -27431.071

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 9990
$ws.Range("J5").Value = 9990
$ws.Range("L5").Value = 9990
$ws.Range("N5").Value = -10214

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 20000
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 20000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 20000
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -20280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1359.1305
$ws.Range("I126").Value = 1330
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 3990
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -1520
$ws.Range("N126").Value = -10940
